$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 277 (shifts the existing rows 277-397 down to 278-398).
$ws.Rows(277).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A277").Value = 6
$ws.Range("B277").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C277").Value = "Metropolitana"
$ws.Range("D277").Value = 44636
$ws.Range("E277").Value = 13
$ws.Range("F277").Value = 100112039
$ws.Range("G277").Value = "Ciboulette"
$ws.Range("H277").Value = "Sin especificar"
$ws.Range("I277").Value = "Primera"
$ws.Range("J277").Value = 620
$ws.Range("K277").Value = 1300
$ws.Range("L277").Value = 1500
$ws.Range("M277").Value = 1406
$ws.Range("N277").Value = "`$/docena de atados"
$ws.Range("O277").Value = "Región Metropolitana"
$ws.Range("P277").Value = 469
$ws.Range("Q277").Value = 3
$ws.Range("R277").Value = "Hortaliza"
